# Applies the cryptos price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.740.19'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '2.463.53'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('E4').Value = '  +0.08%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '573.36'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.65%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '147.81'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.73%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -1.58%  '
$ws.Range('E9').Value = '  +0.31%  '
$ws.Range('E10').Value = '  -0.33%  '
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('E12').Value = '  -0.37%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '29.10'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +1.00%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.0000176'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -1.85%  '
$ws.Range('D16').Value = '62.749.02'
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('D17').Value = '2.467.54'
$ws.Range('E17').Value = '  +0.23%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '7.89'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -1.08%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '10.91'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -1.95%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '325.48'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.41%  '
$ws.Range('E21').Value = '  -0.05%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '2.17'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -2.36%  '
$ws.Range('E23').Value = '  +0.00%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '9.93'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +12.46%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '65.43'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -1.56%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '641.51'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -3.13%  '
$ws.Range('D27').Value = '2.585.94'
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('D28').Value = '0.0₃0974'
$ws.Range('E28').Value = '  -2.86%  '
$ws.Range('E29').Value = '  -13.03%  '
$ws.Range('E30').Value = '  -0.10%  '
$ws.Range('E31').Value = '  -3.12%  '
$ws.Range('E32').Value = '  -3.90%  '
$ws.Range('E33').Value = '  -4.46%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  -0.89%  '
$ws.Range('E36').Value = '  -1.06%  '
$ws.Range('E37').Value = '  -1.53%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '150.79'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -1.48%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '18.59'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -1.42%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '5.31'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -3.97%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.72'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.54%  '
$ws.Range('E42').Value = '  -2.30%  '
$ws.Range('D43').Value = '0.0₆0314'
$ws.Range('E43').Value = '  -11.24%  '
$ws.Range('E44').Value = '  -0.04%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '153.05'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +4.44%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '15.30'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +1.05%  '
$ws.Range('E47').Value = '  -1.54%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '20.34'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -1.84%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.606'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -0.17%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.0508'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -1.64%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0911'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -1.33%  '
